$d = $word.ActiveDocument

# Remove ", application_date" that follows "coupon_id" in the Order tuple.
$d.Content.Find.Execute(", application_date)", $true, $false, $false, $false, $false,
                         $true, 1, $false, ")", 2)
